$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts the existing columns
# (A = labels, B:F = data) one place to the right, becoming B:G.
$ws.Range("A1").EntireColumn.Insert()

# New header cell B1 gets the "segments" label, formatted like the other
# header cells (bold, bordered, centered - same style as C1).
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B1").Value = "segments"

# Fill the new column A (rows 2-20) with a 0-based row index, matching
# the formatting of the label column (now column B).
for ($i = 2; $i -le 20; $i++) {
    $ws.Cells.Item($i, 2).Copy()
    $ws.Cells.Item($i, 1).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Cells.Item($i, 1).Value = $i - 2
}

$excel.CutCopyMode = 0
